$d = $word.ActiveDocument

# Each "Reflective Analysis" bullet currently has its leading label in a
# bold run ("Label: ") followed by a separate, non-bold run with the rest
# of the sentence. Merge them into a single, non-bold run per bullet by
# deleting the bold "Label: " run and re-inserting plain "Label: " text
# in front of the remaining (already non-bold) run.

$items = @(
    @{ label = "Assumptions matter: "; anchor = "Normality and homogeneity" },
    @{ label = "Choice of test: ";      anchor = "Paired vs independent t-tests" },
    @{ label = "Ethical reporting: ";   anchor = "Transparency about significance levels" },
    @{ label = "Research application: "; anchor = "In my SME IS study" }
)

foreach ($item in $items) {
    $delRange = $d.Content
    $found = $delRange.Find.Execute($item.label, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $delRange.Delete()
    }

    $insRange = $d.Content
    $found2 = $insRange.Find.Execute($item.anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $insRange.Collapse(1)
        $insRange.InsertBefore($item.label)
    }
}
